$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert 5 new rows before the old row 5 (pushes everything else
#    down by 5 rows, table grows from 10 to 15 data rows)
# ---------------------------------------------------------------
$ws.Rows("5:9").Insert()

# ---------------------------------------------------------------
# 2. Resize the table (ListObject) to cover the new range
# ---------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F16"))

# ---------------------------------------------------------------
# 3. Update the content of rows 2-4 (identifiers stay the same,
#    datas input cleared, expected output text changed)
# ---------------------------------------------------------------
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "A puse should be sent to both whells motors"
$ws.Rows.Item(2).RowHeight = 30

$ws.Range("E3").Value = "A reverse pulse should be sent to the both whells motors to do a 90°"
$ws.Rows.Item(3).RowHeight = 30

$ws.Range("E4").Value = "A reverse pulse should be sent to the both whells motors to do a 180°"
$ws.Rows.Item(4).RowHeight = 30

# ---------------------------------------------------------------
# 4. Fill in the 5 brand-new rows (5-9)
# ---------------------------------------------------------------
$ws.Range("A5").Value = "DDR_00200"
$ws.Range("B5").Value = "TU_00200"
$ws.Range("C5").Value = "Test if ultrason sensor detects an obstacle"
$ws.Range("E5").Value = "A pulse should be sent to the ultrason sensor when he detects an obstacle to take measure"
$ws.Rows.Item(5).RowHeight = 45

$ws.Range("A6").Value = "DDR_00205"
$ws.Range("B6").Value = "TU_00205"
$ws.Range("C6").Value = "Test if an obstacle is considered detected when the value is under than a threshold (8)"
$ws.Range("D6").Value = "2,`n4,`n8,`n10,`n12,`n14"
$ws.Range("E6").Value = "If input is greater than 8 it should return 0, else it should return 1"
$ws.Rows.Item(6).RowHeight = 90

$ws.Range("A7").Value = "DDR_00210"
$ws.Range("B7").Value = "TU_00210"
$ws.Range("C7").Value = "Test if the servomotor rotates [-45°;+45] clockwise"
$ws.Range("E7").Value = "A pulse should be sent on the servomotor (clockwise pulse)"
$ws.Rows.Item(7).RowHeight = 30

$ws.Range("A8").Value = "DDR_00210"
$ws.Range("B8").Value = "TU_00211"
$ws.Range("C8").Value = "Test if the servomotor rotates [-45°;+45] counterclockwise"
$ws.Range("E8").Value = "A pulse should be sent on the servomotor (counterclockwise pulse)"
$ws.Rows.Item(8).RowHeight = 45

$ws.Range("A9").Value = "DDR_00300"
$ws.Range("B9").Value = "TU_00300"
$ws.Range("C9").Value = "Test if ultrason sensor detects a hole"
$ws.Range("E9").Value = "A pulse should be sent to the infrared sensor  when he detects a hole to take measure"
$ws.Rows.Item(9).RowHeight = 45

# ---------------------------------------------------------------
# 5. Column widths: column C grows a bit wider because of the
#    longer text it now holds
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 33.7

# ---------------------------------------------------------------
# 6. Conditional formatting on the Pass/Fail column: green fill
#    when it contains "PASS", red fill when it contains "FAIL"
# ---------------------------------------------------------------
$rngPF = $ws.Range("F2:F16")

$fcFail = $rngPF.FormatConditions.Add(9, 0, "FAIL")
$fcFail.Text = "FAIL"
$fcFail.Formula1 = 'NOT(ISERROR(SEARCH("FAIL",F2)))'
$fcFail.Interior.Color = 255

$fcPass = $rngPF.FormatConditions.Add(9, 0, "PASS")
$fcPass.Text = "PASS"
$fcPass.Formula1 = 'NOT(ISERROR(SEARCH("PASS",F2)))'
$fcPass.Interior.Color = 5287936

# ---------------------------------------------------------------
# 7. Sheet view: select D6, scroll back to column A
# ---------------------------------------------------------------
$ws.Range("D6").Select()
